$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet holds a daily time series ending at row 139 (date 2025-02-07).
# Append two more days, continuing the series with the same values as the
# last existing row (only column A's date advances).
$srcRow = 139
$newRows = @(140, 141)

foreach ($r in $newRows) {
    $src = $ws.Range("A" + $srcRow + ":J" + $srcRow)
    $dst = $ws.Range("A" + $r + ":J" + $r)
    $src.Copy($dst)

    # Advance the date serial in column A by one day from the row above.
    $ws.Cells.Item($r, 1).Value2 = $ws.Cells.Item($r - 1, 1).Value2 + 1
}
